$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column stores its values as plain text (e.g. "4/9/2025"), not
# as real Excel dates. Force text format on these cells first so Excel's
# COM layer doesn't silently coerce the date-looking strings into date
# serial numbers.
$ws.Range("C2:C5").NumberFormat = "@"

# Row 2: Salary -> Factory Business (amount + date updated)
$ws.Range("A2").Value = "Factory Business"
$ws.Range("B2").Value = 145000
$ws.Range("C2").Value = "4/14/2025"

# Row 3 (new): Taxi Business
$ws.Range("A3").Value = "Taxi Business"
$ws.Range("B3").Value = 98000
$ws.Range("C3").Value = "4/13/2025"

# Row 4 (new): Job Salary
$ws.Range("A4").Value = "Job Salary"
$ws.Range("B4").Value = 50000
$ws.Range("C4").Value = "4/8/2025"

# Row 5: Bank Interest Amount shifted down from row 3 (values unchanged)
$ws.Range("A5").Value = "Bank Interest Amount"
$ws.Range("B5").Value = 15000
$ws.Range("C5").Value = "3/31/2025"
